$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: values become numeric (1025, 1040) instead of text "1025.0"/"1040.0"
$ws.Range("A2").Value = 1025
$ws.Range("B2").Value = 1040

# Row 3: new numeric values (2000, 2000)
$ws.Range("A3").Value = 2000
$ws.Range("B3").Value = 2000

# Row 4: new text values "2000.0" / "2000.0" (stored as text, not numeric,
# leading apostrophe forces text storage instead of number auto-detection)
$ws.Range("A4").Value = "'2000.0"
$ws.Range("B4").Value = "'2000.0"
